$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new (blank) rows at position 12, pushing the existing rows
# 12-96 down to 14-98. This grows the sheet's used range / dimension to
# A1:R98 and renumbers the rest of the weekly price records automatically.
$ws.Rows("12:13").Insert()

# Row 12: new weekly "Primera" quality record (date serial 44532 = 2021-12-02).
$ws.Range("A12").Value = 2
$ws.Range("B12").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44532
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = 100112043
$ws.Range("G12").Value = "Pepino ensalada"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 6500
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 6750
$ws.Range("N12").Value = "$/caja 70 unidades"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 96
$ws.Range("Q12").Value = 70
$ws.Range("R12").Value = "Hortaliza"

# Row 13: new weekly "Segunda" quality record (same date).
$ws.Range("A13").Value = 2
$ws.Range("B13").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44532
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 100112043
$ws.Range("G13").Value = "Pepino ensalada"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Segunda"
$ws.Range("J13").Value = 240
$ws.Range("K13").Value = 4500
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = 4750
$ws.Range("N13").Value = "$/caja 100 unidades"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 48
$ws.Range("Q13").Value = 100
$ws.Range("R13").Value = "Hortaliza"
